$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.204509258270264
$ws.Range("B1").Value = 2.018767356872559
$ws.Range("C1").Value = 4.270934581756592
$ws.Range("D1").Value = 3.015425205230713
$ws.Range("E1").Value = 1.201038718223572
